# Update the cached "average experiment time" style statistics that were
# computed incorrectly (commit: "calc correct avg exp times").
#
# Columns:
#   G = Avg_Agent_Step_Time, H = Avg_Experiment_Time
#   M = Std_Agent_Step_Time, N = Std_Experiment_Time
#
# Only the raw cached numeric values change (no formulas are involved in
# this sheet), so we simply overwrite each cell with the corrected value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G2"  = 6.40288149
    "H2"  = 373.98178118
    "M2"  = 0.5502426597222885
    "N2"  = 55.91611202929847

    "G3"  = 6.78317677
    "H3"  = 610.8127203800001
    "M3"  = 0.6792581408243431
    "N3"  = 114.4994539192471

    "G4"  = 2.53737894
    "H4"  = 71.83620431
    "M4"  = 0.3919814959048554
    "N4"  = 21.28889590035775

    "G5"  = 2.79238079
    "H5"  = 134.33426993
    "M5"  = 0.3614862467633796
    "N5"  = 33.534404364962

    "G6"  = 0.8456443100000001
    "H6"  = 12.93354327
    "M6"  = 0.2371911856015121
    "N6"  = 6.212888359056365

    "G7"  = 1.00181138
    "H7"  = 25.66063576
    "M7"  = 0.1950750737072204
    "N7"  = 8.607412685903103

    "G8"  = 0.41613435
    "H8"  = 4.070704539999999
    "M8"  = 0.1450668131519251
    "N8"  = 2.221703299735914

    "G9"  = 0.48189639
    "H9"  = 8.538840759999999
    "M9"  = 0.1188358353462049
    "N9"  = 3.785142394448388

    "G10" = 0.2117835
    "H10" = 1.53689199
    "M10" = 0.09005140704320205
    "N10" = 0.970451232026567

    "G11" = 0.25371109
    "H11" = 3.546879979999999
    "M11" = 0.07577188831599005
    "N11" = 2.081288556471711

    "G12" = 0.12865917
    "H12" = 0.7916509599999999
    "M12" = 0.06378650906225743
    "N12" = 0.5961162616030753

    "G13" = 0.14540667
    "H13" = 1.63893001
    "M13" = 0.04765479951408701
    "N13" = 0.9543881149747753
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
